$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 150, shifting the existing rows 150-157 down to 151-158
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new record
$ws.Cells.Item(150, 1).Value = 9
$ws.Cells.Item(150, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(150, 3).Value = "Metropolitana"
$ws.Cells.Item(150, 4).Value = 44578
$ws.Cells.Item(150, 5).Value = 13
$ws.Cells.Item(150, 6).Value = "Fruta"
$ws.Cells.Item(150, 7).Value = 100101
$ws.Cells.Item(150, 8).Value = "Berries"
$ws.Cells.Item(150, 9).Value = 100101001
$ws.Cells.Item(150, 10).Value = "Arándano (blue)"
$ws.Cells.Item(150, 11).Value = "Sin especificar"
$ws.Cells.Item(150, 12).Value = "Primera"
$ws.Cells.Item(150, 13).Value = 350
$ws.Cells.Item(150, 14).Value = 4000
$ws.Cells.Item(150, 15).Value = 4000
$ws.Cells.Item(150, 16).Value = 4000
$ws.Cells.Item(150, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(150, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(150, 19).Value = 2000
$ws.Cells.Item(150, 20).Value = 2
